$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to store a literal TEXT value even when the string
    # looks like a number (e.g. "3.04"), matching the original workbook's
    # convention of keeping every data value as a shared string. Toggling
    # the format to Text just long enough to assign the value, then
    # resetting the style back to Normal, avoids leaving any visible
    # number-format / style change behind on the cell itself.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Sheet "Restricciones_del_follower" (3rd sheet) ---------------------
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3.Cells.Item(2,1) "2.9999999999999996 - 2x_1 + y_1 - y_2"
Set-TextValue $ws3.Cells.Item(2,2) "-0.49999999999999956"
Set-TextValue $ws3.Cells.Item(2,4) "0.96"
Set-TextValue $ws3.Cells.Item(2,6) "8.4"

Set-TextValue $ws3.Cells.Item(3,1) "-0.9499999999999997 + x_1 - 3x_2 + y_2"
Set-TextValue $ws3.Cells.Item(3,2) "-1.0500000000000003"
Set-TextValue $ws3.Cells.Item(3,4) "0.9"
Set-TextValue $ws3.Cells.Item(3,5) "3.2"
Set-TextValue $ws3.Cells.Item(3,6) "2.3000000000000003"

Set-TextValue $ws3.Cells.Item(4,1) "-4.14 + x_1 + x_2"
Set-TextValue $ws3.Cells.Item(4,2) "1.7999999999999998"
Set-TextValue $ws3.Cells.Item(4,4) "0.28"
Set-TextValue $ws3.Cells.Item(4,5) "6.7"
Set-TextValue $ws3.Cells.Item(4,6) "7.1"

# --- Sheet "Punto_modificado" (4th sheet) --------------------------------
$ws4 = $wb.Worksheets.Item(4)

Set-TextValue $ws4.Cells.Item(2,1) "2.55"
Set-TextValue $ws4.Cells.Item(2,2) "1.25"
Set-TextValue $ws4.Cells.Item(2,3) "4.25"
Set-TextValue $ws4.Cells.Item(2,4) "2.15"

# --- Sheet "Vector_bf" (5th sheet) ---------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-TextValue $ws5.Cells.Item(2,1) "3.04"
Set-TextValue $ws5.Cells.Item(3,1) "-0.9400000000000001"

# --- Sheet "Vector_BF" (6th sheet) ---------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-TextValue $ws6.Cells.Item(2,1) "-7.9"
Set-TextValue $ws6.Cells.Item(3,1) "1.9000000000000012"
Set-TextValue $ws6.Cells.Item(5,1) "-3.2"
